# Phase-1 Commit and Push
# Update the Ticket Number column (Y) on the NFTRTickets sheet with the
# newly generated 2024 ticket numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFTRTickets")

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "240720001008"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "240720001009"

$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "240720001010"

$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "240720001011"
